# Applies the "annual electricity consumption reduced to 1200+" edit:
#  - Updates several C-column (value) cells for ID_Technology rows 22-28 and 35-36
#  - Removes the last two data rows (ID_Technology 38 and 39)
#  - Restores a clean alternating row-shading pattern on column A starting at row 23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update changed values in column C (row -> new value) ---
$valueChanges = @{
    23 = 14
    24 = 1
    25 = 0.3
    26 = 0.3
    28 = 1
    29 = 1
    36 = 0
    37 = 0
}

foreach ($row in $valueChanges.Keys) {
    $ws.Cells.Item($row, 3).Value2 = $valueChanges[$row]
}

# --- 2) Delete the trailing two rows (old rows 39 and 40) ---
$ws.Rows.Item(40).Delete() | Out-Null
$ws.Rows.Item(39).Delete() | Out-Null

# --- 3) Reapply alternating shading to column A for rows 23-38 ---
# Style "2": plain font, no fill.  Style "3": plain font, light-grey fill (RGB E7E6E6).
$shadeColor = 15132391   # OLE (BGR) value of RGB FFE7E6E6

for ($row = 23; $row -le 38; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Font.Color = 0
    $shouldShade = ((($row - 23) % 2) -eq 0)
    if ($shouldShade) {
        $cell.Interior.Color = $shadeColor
        $cell.Interior.Pattern = 1
    } elseif ($cell.Interior.Pattern -ne -4142) {
        # only needs a reset when the cell currently carries a fill that must be cleared
        $cell.Interior.ColorIndex = -4142
        $cell.Interior.Pattern = -4142
    }
}

# --- 4) Restore the cursor/selection to E15, matching the saved view state ---
$ws.Range("E15").Select() | Out-Null
